# Split of File Access Screen file into two files.
# One for accessing the files and another to grade.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 5 (aid 7 -> 24, reassigned to a different access record) ---
$ws.Range("A5").Value = 24
$ws.Range("E5").Value = 16208102
$ws.Range("F5").Value = "test1.txt"
$ws.Range("G5").Value = 40
$ws.Range("I5").Value = 43936.750031970783

# --- Update existing row 13 (aid 5 -> 26, reassigned to a different access record) ---
$ws.Range("A13").Value = 26
$ws.Range("E13").Value = 111111111
$ws.Range("F13").Value = "test2.txt"
$ws.Range("I13").Value = 43937.96865753649

# --- Add new row 14: the original aid 7 grading record that was split out ---
$ws.Range("A14").Value = 7
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "CS4125"
$ws.Range("D14").Value = 1
$ws.Range("E14").Value = 12345678
$ws.Range("F14").Value = "DocTest.txt"
$ws.Range("G14").Value = 25
$ws.Range("H14").Value = "Y"
$ws.Range("I14").Value = 43932.897177812352

# --- Column width / best-fit cosmetics to mirror the authored workbook ---
$ws.Columns.Item(7).ColumnWidth = 10.14
$ws.Columns.Item(9).ColumnWidth = 18.71

# --- Restore the last active selection recorded in the sheet view ---
$ws.Range("H17").Select() | Out-Null
